$d = $word.ActiveDocument

$rng = $d.Range(277, 277)
$rng.Select()
$sel = $word.Selection
Write-Output ("Selection start/end: " + $sel.Start + "-" + $sel.End)
$sel.TypeText("XYZ")
Write-Output ("Doc length: " + $d.Content.End)
